$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product backlog")

$ws.Range("B4").Value = "Als beheerder wil ik dat de boerderij zichzelf stabiliseert door water te verplaatsen met pompen, zodat dit duurzaam en efficient gebeurd. "
$ws.Range("B3").Value = "Als systeem wil ik de hoek nauwkeurig berekenen, zodat ik weet naar welke kant water verplaatst moet worden."
$ws.Range("B2").Value = "Ik wil duidelijk kunnen zien in welke hoek de boerderij zich bevindt en in welke toestand de andere onderdelen verkeren, om een duidelijk overzicht te hebben van het systeem. "
$ws.Range("B5").Value = "Ik wil graag een duidelijk en werkend schaalmodel zien, waarin goed wordt weergegeven wat er gebeurt wanneer de boerderij uit balans raakt. "
$ws.Range("A6").Value = "5. Algoritme om floating farm waterpas houden"
$ws.Range("B6").Value = "Als systeem wil ik sensordata waterpas hebben, zodat ik altijd nauwkeurige en realistische lezingen krijg."
